$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix month/year column swap: A should hold month number, B should hold year (2020)
for ($r = 2; $r -le 13; $r++) {
    $month = $ws.Cells.Item($r, 2).Value2
    $year = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value2 = $month
    $ws.Cells.Item($r, 2).Value2 = $year
}

# New headers for columns H:P
$ws.Range("H1").Value2 = "grade_total"
$ws.Range("I1").Value2 = "grade_distance"
$ws.Range("J1").Value2 = "grade_visitation"
$ws.Range("K1").Value2 = "grade_encounters"
$ws.Range("L1").Value2 = "NEVER"
$ws.Range("M1").Value2 = "RARELY"
$ws.Range("N1").Value2 = "SOMETIMES"
$ws.Range("O1").Value2 = "FREQUENTLY"
$ws.Range("P1").Value2 = "ALWAYS"

# New data values for columns H:P, rows 2-13 (same values for every row)
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 8).Value2 = 0.7
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = 0
    $ws.Cells.Item($r, 11).Value2 = 1
    $ws.Cells.Item($r, 12).Value2 = 1.023
    $ws.Cells.Item($r, 13).Value2 = 1.024
    $ws.Cells.Item($r, 14).Value2 = 1.073
    $ws.Cells.Item($r, 15).Value2 = 1.121
    $ws.Cells.Item($r, 16).Value2 = 1.759
}
